# Insert a new weekly record as row 498 in the "Pepino ensalada" price sheet,
# pushing the existing rows 498:534 down to 499:535 (dimension grows to R535).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("498:498").Insert()

$ws.Range("A498").Value = 3
$ws.Range("B498").Value = "Femacal de La Calera"
$ws.Range("C498").Value = "Coquimbo"
$ws.Range("D498").Value = 45021
$ws.Range("E498").Value = 5
$ws.Range("F498").Value = 100112043
$ws.Range("G498").Value = "Pepino ensalada"
$ws.Range("H498").Value = "Sin especificar"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 135
$ws.Range("K498").Value = 8500
$ws.Range("L498").Value = 9000
$ws.Range("M498").Value = 8759
$ws.Range("N498").Value = "$/caja 60 unidades"
$ws.Range("O498").Value = "Región de Arica y Parinacota"
$ws.Range("P498").Value = 146
$ws.Range("Q498").Value = 60
$ws.Range("R498").Value = "Hortaliza"
